$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Hombres" / "Mujeres" header text between E1 and F1
$ws.Range("E1").Value = "Mujeres"
$ws.Range("F1").Value = "Hombres"

# Update the active cell selection on the frozen (bottom-left) pane from D4 to F8
$ws.Range("F8").Select()
